# Auto-generated edit script: updates currentAveragePrice / Leve price / profit
# columns (H-N) for the specified rows across multiple sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 2212.5
$ws.Range("I2").Value = 2641.6667
$ws.Range("J2").Value = 925
$ws.Range("K2").Value = 2641.6667
$ws.Range("L2").Value = 925
$ws.Range("M2").Value = -2528.6667
$ws.Range("N2").Value = -1151
# Row 58
$ws.Range("H58").Value = 681.1667
$ws.Range("I58").Value = 217.4
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 652.2
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -502.2
$ws.Range("N58").Value = -9300
# Row 74
$ws.Range("H74").Value = 5421.2
$ws.Range("I74").Value = 5165
$ws.Range("J74").Value = 5714
$ws.Range("K74").Value = 5165
$ws.Range("L74").Value = 5714
$ws.Range("M74").Value = -4229
$ws.Range("N74").Value = -7586
# Row 77
$ws.Range("H77").Value = 5421.2
$ws.Range("I77").Value = 5165
$ws.Range("J77").Value = 5714
$ws.Range("K77").Value = 25825
$ws.Range("L77").Value = 28570
$ws.Range("M77").Value = -21145
$ws.Range("N77").Value = -37930
# Row 116
$ws.Range("H116").Value = 6470.25
$ws.Range("I116").Value = 8000.5
$ws.Range("J116").Value = 4940
$ws.Range("K116").Value = 8000.5
$ws.Range("L116").Value = 4940
$ws.Range("M116").Value = -4558.5
$ws.Range("N116").Value = -11824
# Row 137
$ws.Range("H137").Value = 6885.857
$ws.Range("I137").Value = 15999.75
$ws.Range("J137").Value = 3240.3
$ws.Range("K137").Value = 47999.25
$ws.Range("L137").Value = 9720.900000000001
$ws.Range("M137").Value = -45449.25
$ws.Range("N137").Value = -14820.9

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1348047.6
$ws.Range("I32").Value = 1370105.8
$ws.Range("J32").Value = 2500
$ws.Range("K32").Value = 1370105.8
$ws.Range("L32").Value = 2500
$ws.Range("M32").Value = -1369818.8
$ws.Range("N32").Value = -3074
# Row 97
$ws.Range("H97").Value = 647.5714
$ws.Range("I97").Value = 670.5
$ws.Range("J97").Value = 349.5
$ws.Range("K97").Value = 670.5
$ws.Range("L97").Value = 349.5
$ws.Range("M97").Value = -174.5
$ws.Range("N97").Value = -1341.5
# Row 122
$ws.Range("H122").Value = 2586
$ws.Range("I122").Value = 2521.3
$ws.Range("J122").Value = 2801.6667
$ws.Range("K122").Value = 7563.900000000001
$ws.Range("L122").Value = 8405.000100000001
$ws.Range("M122").Value = -5113.900000000001
$ws.Range("N122").Value = -13305.0001

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 12230.6
$ws.Range("I86").Value = 50000
$ws.Range("J86").Value = 2788.25
$ws.Range("K86").Value = 50000
$ws.Range("L86").Value = 2788.25
$ws.Range("M86").Value = -48877
$ws.Range("N86").Value = -5034.25
# Row 89
$ws.Range("H89").Value = 12230.6
$ws.Range("I89").Value = 50000
$ws.Range("J89").Value = 2788.25
$ws.Range("K89").Value = 250000
$ws.Range("L89").Value = 13941.25
$ws.Range("M89").Value = -244384
$ws.Range("N89").Value = -25173.25
# Row 94
$ws.Range("H94").Value = 947.53845
$ws.Range("I94").Value = 659.7778
$ws.Range("J94").Value = 1595
$ws.Range("K94").Value = 659.7778
$ws.Range("L94").Value = 1595
$ws.Range("M94").Value = -208.7778
$ws.Range("N94").Value = -2497
# Row 100
$ws.Range("H100").Value = 16160
$ws.Range("J100").Value = 16160
$ws.Range("L100").Value = 16160
$ws.Range("N100").Value = -18324

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4996.7144
$ws.Range("I31").Value = 1779.9445
$ws.Range("J31").Value = 6864.516
$ws.Range("K31").Value = 1779.9445
$ws.Range("L31").Value = 6864.516
$ws.Range("M31").Value = -1484.9445
$ws.Range("N31").Value = -7454.516
# Row 34
$ws.Range("H34").Value = 4996.7144
$ws.Range("I34").Value = 1779.9445
$ws.Range("J34").Value = 6864.516
$ws.Range("K34").Value = 1779.9445
$ws.Range("L34").Value = 6864.516
$ws.Range("M34").Value = -1577.9445
$ws.Range("N34").Value = -7268.516
# Row 43
$ws.Range("H43").Value = 14950
$ws.Range("J43").Value = 14950
$ws.Range("L43").Value = 14950
$ws.Range("N43").Value = -15318
# Row 101
$ws.Range("H101").Value = 14950
$ws.Range("J101").Value = 14950
$ws.Range("L101").Value = 14950
$ws.Range("N101").Value = -21440

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1407.2084
$ws.Range("I5").Value = 861.38464
$ws.Range("J5").Value = 2052.2727
$ws.Range("K5").Value = 2584.15392
$ws.Range("L5").Value = 6156.8181
$ws.Range("M5").Value = -2472.15392
$ws.Range("N5").Value = -6380.8181
# Row 107
$ws.Range("H107").Value = 1308.25
$ws.Range("I107").Value = 210.75
$ws.Range("J107").Value = 1674.0834
$ws.Range("K107").Value = 632.25
$ws.Range("L107").Value = 5022.2502
$ws.Range("M107").Value = 1287.75
$ws.Range("N107").Value = -8862.2502
# Row 135
$ws.Range("H135").Value = 1407.2084
$ws.Range("I135").Value = 861.38464
$ws.Range("J135").Value = 2052.2727
$ws.Range("K135").Value = 7752.46176
$ws.Range("L135").Value = 18470.4543
$ws.Range("M135").Value = -5217.46176
$ws.Range("N135").Value = -23540.4543

$ws = $wb.Worksheets.Item("GSM")
# Row 98
$ws.Range("H98").Value = 25079.9
$ws.Range("J98").Value = 25079.9
$ws.Range("L98").Value = 25079.9
$ws.Range("N98").Value = -31069.9

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 640.4
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 701
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 701
$ws.Range("M22").Value = -305
$ws.Range("N22").Value = -1291
# Row 27
$ws.Range("H27").Value = 640.4
$ws.Range("I27").Value = 600
$ws.Range("J27").Value = 701
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 701
$ws.Range("M27").Value = -493
$ws.Range("N27").Value = -915
# Row 122
$ws.Range("H122").Value = 1826.2
$ws.Range("I122").Value = 1863
$ws.Range("J122").Value = 1784.1428
$ws.Range("K122").Value = 5589
$ws.Range("L122").Value = 5352.428400000001
$ws.Range("M122").Value = -3139
$ws.Range("N122").Value = -10252.4284

$ws = $wb.Worksheets.Item("WVR")
# Row 63
$ws.Range("H63").Value = 21721
$ws.Range("J63").Value = 21721
$ws.Range("L63").Value = 21721
$ws.Range("N63").Value = -22969
# Row 66
$ws.Range("H66").Value = 21721
$ws.Range("J66").Value = 21721
$ws.Range("L66").Value = 65163
$ws.Range("N66").Value = -71403
# Row 92
$ws.Range("H92").Value = 30275
$ws.Range("J92").Value = 30275
$ws.Range("L92").Value = 30275
$ws.Range("N92").Value = -35267
# Row 93
$ws.Range("H93").Value = 35000
$ws.Range("J93").Value = 35000
$ws.Range("L93").Value = 35000
$ws.Range("N93").Value = -39992
# Row 96
$ws.Range("H96").Value = 2355.0588
$ws.Range("I96").Value = 2411
$ws.Range("J96").Value = 2292.125
$ws.Range("K96").Value = 2411
$ws.Range("L96").Value = 2292.125
$ws.Range("M96").Value = -1038
$ws.Range("N96").Value = -5038.125
# Row 99
$ws.Range("H99").Value = 22000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 22000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 22000
$ws.Range("M99").Value = ""
$ws.Range("N99").Value = -27990
# Row 105
$ws.Range("H105").Value = 44629.332
$ws.Range("J105").Value = 44629.332
$ws.Range("L105").Value = 44629.332
$ws.Range("N105").Value = -51617.332
# Row 132
$ws.Range("H132").Value = 3278.8
$ws.Range("I132").Value = 2629.8462
$ws.Range("J132").Value = 4484
$ws.Range("K132").Value = 7889.5386
$ws.Range("L132").Value = 13452
$ws.Range("M132").Value = -5359.5386
$ws.Range("N132").Value = -18512

